$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-25 Saturday" "2024-05-26 Sunday"
Replace-Text "119×6=" "776×3="
Replace-Text "173×8=" "589×3="
Replace-Text "408×2=" "567×4="
Replace-Text "501×3=" "749×7="
Replace-Text "715×7=" "750×3="
Replace-Text "241×8=" "770×6="
Replace-Text "614×3=" "715×9="
Replace-Text "480×9=" "358×9="
Replace-Text "655×5=" "876×9="
Replace-Text "621×8=" "517×6="
Replace-Text "980×9=" "733×5="
Replace-Text "692×3=" "510×3="
Replace-Text "579×8=" "579×7="
Replace-Text "553×5=" "109×9="
Replace-Text "302×3=" "483×6="
Replace-Text "872×7=" "705×9="
Replace-Text "831×8=" "336×5="
Replace-Text "207×8=" "707×9="
Replace-Text "896×6=" "236×6="
Replace-Text "521×7=" "215×4="
Replace-Text "439×7=" "117×2="
Replace-Text "377×9=" "242×5="
Replace-Text "766×4=" "620×2="
Replace-Text "622×9=" "443×5="
Replace-Text "272×5=" "246×4="
